$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C49").Value = "30 - 40 lacs per annum"
$ws.Range("C50").Value = "50 - 80 lacs per annum"
$ws.Range("C51").Value = "> 80 lacs per annum"
$ws.Range("D51").Value = "why care?"
$ws.Range("D49").Value = "6 - 10 lacs"
$ws.Range("D50").Value = "12 - 36 lacs"
